# Append a new "2021年" data row (row 11) to Sheet1, right after the
# existing last row (row 10, "2020年"), keeping the same A:AQ layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 11

# Column A: year label, styled the same way as the other year cells (A2:A10).
$ws.Range("A" + ($newRow - 1)).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item($newRow, 1).Value = "2021年"

# Numeric values for columns B..AQ. Columns E (5) and V (22) are left blank,
# matching every other data row in the sheet.
$rowData = [ordered]@{
    2  = 10380.78            # B
    3  = 2524.87             # C
    4  = 428.56              # D
    5  = $null                # E (blank)
    6  = 6246.66              # F
    7  = 18562.59             # G
    8  = 1944.42              # H
    9  = 10887.88             # I
    10 = 1504.11              # J
    11 = 279178.61            # K
    12 = 1411.57              # L
    13 = 313.92               # M
    14 = 90.84                # N
    15 = 2729.9               # O
    16 = 5079.2               # P
    17 = 366.92               # Q
    18 = 429.64               # R
    19 = 6930.6               # S
    20 = 2579.2               # T
    21 = 37713.49             # U
    22 = $null                 # V (blank)
    23 = 3307.72              # W
    24 = 5036.01              # X
    25 = 13124.16             # Y
    26 = 17332                # Z
    27 = 2138.65              # AA
    28 = 3740.34              # AB
    29 = 3246.26              # AC
    30 = 4019.41              # AD
    31 = 3464.67              # AE
    32 = 59174.41             # AF
    33 = 13975.02             # AG
    34 = 6448.31              # AH
    35 = 3303.05              # AI
    36 = 565.9400000000001    # AJ
    37 = 7290.48              # AK
    38 = 3373.46              # AL
    39 = 7650.86              # AM
    40 = 136.83               # AN
    41 = 6111.96              # AO
    42 = 5340.73              # AP
    43 = 237.9                # AQ
}

foreach ($col in $rowData.Keys) {
    $cell = $ws.Cells.Item($newRow, $col)
    $value = $rowData[$col]
    if ($null -eq $value) {
        # Keep the cell present (like the other blank cells in this column)
        # but with no value, instead of leaving it absent entirely.
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
